$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking text values (Price / Volume columns) stay as text
$textCells = @('D2', 'E2', 'D3', 'E3', 'D4', 'E4', 'D5', 'E5', 'D6', 'E6', 'D7', 'E7', 'D8', 'E8', 'D9', 'E9', 'D10', 'E10', 'D11', 'E11', 'D12', 'E12', 'E13', 'D14', 'E14', 'D15', 'E15', 'D16', 'E16', 'D17', 'E17', 'D18', 'E18', 'D19', 'E19', 'D20', 'E20', 'D21', 'E21', 'D22', 'E22', 'D23', 'E23', 'D24', 'E24', 'E25', 'D26', 'D38', 'E38', 'D39', 'E39', 'E40', 'D41', 'E41', 'D42', 'E42', 'E43', 'D44', 'E44', 'D45', 'E45', 'D46', 'E46', 'E47', 'E49', 'E50', 'E51')
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply updated values
$ws.Range('D2').Value = '289.96'
$ws.Range('E2').Value = '-9.49%'
$ws.Range('D3').Value = '40.19'
$ws.Range('E3').Value = '-2.96%'
$ws.Range('D4').Value = '5.042'
$ws.Range('E4').Value = '-4.11%'
$ws.Range('D5').Value = '0.07306'
$ws.Range('E5').Value = '-5.66%'
$ws.Range('D6').Value = '4.276'
$ws.Range('E6').Value = '-1.37%'
$ws.Range('D7').Value = '1.557'
$ws.Range('E7').Value = '-11.45%'
$ws.Range('D8').Value = '0.9182'
$ws.Range('E8').Value = '-2.81%'
$ws.Range('D9').Value = '0.1164'
$ws.Range('E9').Value = '-7.68%'
$ws.Range('D10').Value = '0.1721'
$ws.Range('E10').Value = '-7.66%'
$ws.Range('D11').Value = '0.08692'
$ws.Range('E11').Value = '-5.64%'
$ws.Range('D12').Value = '0.04164'
$ws.Range('E12').Value = '0.31%'
$ws.Range('E13').Value = '0.28%'
$ws.Range('D14').Value = '0.001261'
$ws.Range('E14').Value = '-1.73%'
$ws.Range('D15').Value = '0.005784'
$ws.Range('E15').Value = '-0.70%'
$ws.Range('B16').Value = 'LEO'
$ws.Range('C16').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D16').Value = '3.393'
$ws.Range('E16').Value = '1.24%'
$ws.Range('B17').Value = 'BTSEToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D17').Value = '2.397'
$ws.Range('E17').Value = '-1.16%'
$ws.Range('B18').Value = 'BitpandaEcosystemToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range('D18').Value = '0.3277'
$ws.Range('E18').Value = '-2.40%'
$ws.Range('B19').Value = 'MCDex'
$ws.Range('C19').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('D19').Value = '7.801'
$ws.Range('E19').Value = '-7.30%'
$ws.Range('B20').Value = 'ProBitToken'
$ws.Range('C20').Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range('D20').Value = '0.1350'
$ws.Range('E20').Value = '-0.24%'
$ws.Range('B21').Value = 'ZBToken'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range('D21').Value = '0.2882'
$ws.Range('E21').Value = '2.00%'
$ws.Range('B22').Value = 'CoinExToken'
$ws.Range('C22').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D22').Value = '0.03863'
$ws.Range('E22').Value = '-4.28%'
$ws.Range('B23').Value = 'BitKan'
$ws.Range('C23').Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range('D23').Value = '0.001268'
$ws.Range('E23').Value = '-0.09%'
$ws.Range('B24').Value = 'HotbitToken'
$ws.Range('C24').Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range('D24').Value = '0.003890'
$ws.Range('E24').Value = '-5.57%'
$ws.Range('E25').Value = '0.67%'
$ws.Range('D26').Value = '0.0003723'
$ws.Range('D38').Value = '0.02327'
$ws.Range('E38').Value = '-8.83%'
$ws.Range('D39').Value = '0.04964'
$ws.Range('E39').Value = '-7.30%'
$ws.Range('E40').Value = '237.38%'
$ws.Range('D41').Value = '0.007706'
$ws.Range('E41').Value = '-0.88%'
$ws.Range('D42').Value = '0.1276'
$ws.Range('E42').Value = '-3.24%'
$ws.Range('E43').Value = '4.86%'
$ws.Range('D44').Value = '0.007064'
$ws.Range('E44').Value = '-14.95%'
$ws.Range('D45').Value = '0.2891'
$ws.Range('E45').Value = '-16.26%'
$ws.Range('D46').Value = '0.00006414'
$ws.Range('E46').Value = '-4.14%'
$ws.Range('E47').Value = '-0.06%'
$ws.Range('E49').Value = '-85.02%'
$ws.Range('E50').Value = '-0.06%'
$ws.Range('E51').Value = '-0.06%'
